$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force cells that would otherwise be auto-parsed as numbers to stay as text,
# matching the original inline-string (text) representation in the source data.
$textCells = @("D5","D6","D8","D12","D13","D16","D18","D19","D20","D23","D25","D26","D28","D30","D33","D35","D36","D37","D42","D43","D45","D46","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "67.764.22"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "2.498.94"
$ws.Range("E3").Value = "  +0.92%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "588.15"
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("D6").Value = "176.76"
$ws.Range("E6").Value = "  +4.46%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.516"
$ws.Range("E8").Value = "  +1.01%  "
$ws.Range("E9").Value = "  +6.05%  "
$ws.Range("E10").Value = "  +0.40%  "
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").Value = "4.95"
$ws.Range("E12").Value = "  +0.81%  "
$ws.Range("D13").Value = "25.81"
$ws.Range("E13").Value = "  +2.21%  "
$ws.Range("D14").Value = "2.917.22"
$ws.Range("E14").Value = "  -0.09%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "67.635.89"
$ws.Range("E15").Value = "  +1.07%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("E16").Value = "  +2.56%  "
$ws.Range("D17").Value = "2.528.68"
$ws.Range("E17").Value = "  +2.81%  "
$ws.Range("D18").Value = "11.11"
$ws.Range("E18").Value = "  +1.62%  "
$ws.Range("D19").Value = "7.55"
$ws.Range("E19").Value = "  +2.71%  "
$ws.Range("D20").Value = "351.72"
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  +2.24%  "
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "70.69"
$ws.Range("E23").Value = "  +3.28%  "
$ws.Range("E24").Value = "  +3.01%  "
$ws.Range("D25").Value = "1.77"
$ws.Range("E25").Value = "  -1.54%  "
$ws.Range("D26").Value = "9.18"
$ws.Range("E26").Value = "  +0.75%  "
$ws.Range("E27").Value = "  +1.42%  "
$ws.Range("D28").Value = "1.00"
$ws.Range("E28").Value = "  +0.44%  "
$ws.Range("D29").Value = "0.0₃0914"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").Value = "509.05"
$ws.Range("E30").Value = "  +0.04%  "
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("E32").Value = "  +3.26%  "
$ws.Range("D33").Value = "1.78"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  +0.09%  "
$ws.Range("D35").Value = "0.123"
$ws.Range("E35").Value = "  +7.77%  "
$ws.Range("D36").Value = "163.52"
$ws.Range("E36").Value = "  +3.07%  "
$ws.Range("D37").Value = "18.45"
$ws.Range("E37").Value = "  +1.54%  "
$ws.Range("E38").Value = "  +0.11%  "
$ws.Range("E39").Value = "  +1.72%  "
$ws.Range("E40").Value = "  +0.09%  "
$ws.Range("E41").Value = "  +4.37%  "
$ws.Range("D42").Value = "0.330"
$ws.Range("E42").Value = "  +1.52%  "
$ws.Range("D43").Value = "4.88"
$ws.Range("E43").Value = "  +1.80%  "
$ws.Range("E44").Value = "  +3.12%  "
$ws.Range("D45").Value = "146.14"
$ws.Range("E45").Value = "  +3.63%  "
$ws.Range("D46").Value = "3.52"
$ws.Range("E46").Value = "  +2.80%  "
$ws.Range("E47").Value = "  +1.36%  "
$ws.Range("D48").Value = "0.0₆0257"
$ws.Range("E48").Value = "  +2.90%  "
$ws.Range("D49").Value = "0.0746"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +2.44%  "
$ws.Range("E51").Value = "  +1.11%  "
